{"js": "// Replace the multiplication-problem text in the table cells with the\n// newly generated set of problems. Each old string occurs exactly once\n// in the document, so a targeted search + replace per pair is safe and\n// unambiguous.\nconst replacements = [\n  [\"782\u00d78=6256\", \"525\u00d77=3675\"],\n  [\"542\u00d74=2168\", \"924\u00d77=6468\"],\n  [\"775\u00d79=6975\", \"334\u00d75=1670\"],\n  [\"988\u00d78=7904\", \"429\u00d74=1716\"],\n  [\"543\u00d72=1086\", \"399\u00d76=2394\"],\n  [\"292\u00d79=2628\", \"153\u00d75=765\"],\n  [\"489\u00d76=2934\", \"767\u00d78=6136\"],\n  [\"475\u00d79=4275\", \"857\u00d76=5142\"],\n  [\"935\u00d79=8415\", \"885\u00d77=6195\"],\n  [\"605\u00d77=4235\", \"239\u00d78=1912\"],\n  [\"861\u00d79=7749\", \"272\u00d72=544\"],\n  [\"538\u00d74=2152\", \"574\u00d73=1722\"],\n  [\"627\u00d72=1254\", \"972\u00d76=5832\"],\n  [\"472\u00d75=2360\", \"257\u00d75=1285\"],\n  [\"919\u00d74=3676\", \"742\u00d79=6678\"],\n  [\"849\u00d77=5943\", \"418\u00d76=2508\"],\n  [\"278\u00d77=1946\", \"812\u00d79=7308\"],\n  [\"509\u00d73=1527\", \"946\u00d78=7568\"],\n  [\"610\u00d73=1830\", \"990\u00d76=5940\"],\n  [\"178\u00d77=1246\", \"863\u00d74=3452\"],\n  [\"536\u00d79=4824\", \"927\u00d77=6489\"],\n  [\"560\u00d76=3360\", \"160\u00d73=480\"],\n  [\"636\u00d74=2544\", \"419\u00d73=1257\"],\n  [\"266\u00d72=532\", \"642\u00d79=5778\"],\n  [\"656\u00d79=5904\", \"245\u00d74=980\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const found = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the multiplication-problem text in the table cells with the\n# newly generated set of problems. Each old string occurs exactly once\n# in the document, so a targeted Find/Replace per pair is safe and\n# unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"782\u00d78=6256\", \"525\u00d77=3675\"),\n    @(\"542\u00d74=2168\", \"924\u00d77=6468\"),\n    @(\"775\u00d79=6975\", \"334\u00d75=1670\"),\n    @(\"988\u00d78=7904\", \"429\u00d74=1716\"),\n    @(\"543\u00d72=1086\", \"399\u00d76=2394\"),\n    @(\"292\u00d79=2628\", \"153\u00d75=765\"),\n    @(\"489\u00d76=2934\", \"767\u00d78=6136\"),\n    @(\"475\u00d79=4275\", \"857\u00d76=5142\"),\n    @(\"935\u00d79=8415\", \"885\u00d77=6195\"),\n    @(\"605\u00d77=4235\", \"239\u00d78=1912\"),\n    @(\"861\u00d79=7749\", \"272\u00d72=544\"),\n    @(\"538\u00d74=2152\", \"574\u00d73=1722\"),\n    @(\"627\u00d72=1254\", \"972\u00d76=5832\"),\n    @(\"472\u00d75=2360\", \"257\u00d75=1285\"),\n    @(\"919\u00d74=3676\", \"742\u00d79=6678\"),\n    @(\"849\u00d77=5943\", \"418\u00d76=2508\"),\n    @(\"278\u00d77=1946\", \"812\u00d79=7308\"),\n    @(\"509\u00d73=1527\", \"946\u00d78=7568\"),\n    @(\"610\u00d73=1830\", \"990\u00d76=5940\"),\n    @(\"178\u00d77=1246\", \"863\u00d74=3452\"),\n    @(\"536\u00d79=4824\", \"927\u00d77=6489\"),\n    @(\"560\u00d76=3360\", \"160\u00d73=480\"),\n    @(\"636\u00d74=2544\", \"419\u00d73=1257\"),\n    @(\"266\u00d72=532\", \"642\u00d79=5778\"),\n    @(\"656\u00d79=5904\", \"245\u00d74=980\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    # 2 == wdReplaceAll\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
